# RPA datasets push 2024-06-01
# A new IPO entry ("아이빔테크놀로지") is inserted at the top of the listing
# (row 2), pushing the existing rows down by one and dropping the oldest
# entry that previously sat in the last row (21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row right below the header and strip any formatting it
# may have inherited so it matches the plain (unstyled) data rows below it.
$ws.Rows("2:2").Insert()
$ws.Rows("2:2").ClearFormats()

$ws.Range("A2").Value = "아이빔테크놀로지"
$ws.Range("B2").Value = "2024.07.15~07.19"
$ws.Range("C2").Value = "7,300~8,500"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 16308
$ws.Range("F2").Value = "삼성증권"

# The table keeps a fixed 20-row window, so drop the row that fell off the
# bottom when everything else shifted down.
$ws.Rows("22:22").Delete()
